# prod smoke test debugging - refresh MSRP / DPHF pricing figures on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - GX 460 (9700): BASE MSRP 53000 -> 53100
$ws.Range("D29").Value = 53100

# Row 30 - GX 460 Premium (9700PM): BASE MSRP 55790 -> 55890
$ws.Range("D30").Value = 55890

# Row 31 - GX 460 Luxury (9710): BASE MSRP 64265 -> 64365
$ws.Range("D31").Value = 64365

# Row 32 - LX 570 Two-Row (9625): BASE MSRP 86480 -> 86580, DPHF 1295 -> 1025
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

# Row 33 - LX 570 Three-Row (9620): BASE MSRP 91480 -> 91580, DPHF 1295 -> 1025
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# Row 34 - LX 570 Inspiration Series SE (9620 (SE)): BASE MSRP was blank ("  ")
# -> now a real number 99310, matching the same "#,##0" format used by the
# other BASE MSRP cells in this column; DPHF 1295 -> 1025
$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("D34").Value = 99310
$ws.Range("E34").Value = 1025

# Reflect the author's last-used cell/selection in the saved view state
$ws.Range("D29").Select()
